$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("otv_grubu_co2_araliklari")
$ws1.Range("F2").Value = 2200
$ws1.Range("G2").Formula = "=200+F2"
$ws1.Range("H2").Formula = "=200+G2"
$ws1.Range("I2").Formula = "=200+H2"
$ws1.Range("J2").Formula = "=200+I2"

Write-Output "done"
